$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-7 from 2023-09-16 (45185)
# to 2023-10-05 (45204), keeping the existing date style/format intact.
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
